# Update marksheet totals: correct marks / total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - raw score contribution changes from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row - total marks obtained changes from 66 to 110
$ws.Range("B12").Value = 110

# "Total" row - Correct/Total display string changes from 66/84 to 110/140
$ws.Range("E12").Value = "110/140"
